{"js": "// Mark two maintenance-request items as done, and add \"Done: \" prefixes:\n//   - \"Add custom field size\"          -> \"Done: Add custom field size\"\n//   - \"Randomize \u2013 blocks / num of blocks\" -> \"Done: Randomize \u2013 blocks / num of blocks\"\n// (matches commit message: \"Made change to randomize both the human and zombie\")\n//\n// The new \"Done: \" text must land in its OWN run (not merged into the\n// existing run), so we splice in a tiny OOXML fragment (flat-OPC form,\n// as required by Range.insertOoxml) right before the paragraph's\n// existing content instead of using insertText/insertParagraph.\n\nfunction doneRunOoxml(text) {\n  const escaped = text\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p><w:r>\" +\n    '<w:t xml:space=\"preserve\">' + escaped + \"</w:t>\" +\n    \"</w:r></w:p></w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet addCustomFieldPara = null;\nlet zombieHumanPara = null;\nlet blocksPara = null;\n\nfor (const p of paragraphs.items) {\n  const t = p.text.trim();\n  if (t === \"Add custom field size\") {\n    addCustomFieldPara = p;\n  } else if (t.indexOf(\"Randomize\") === 0 && t.indexOf(\"zombie\") !== -1) {\n    zombieHumanPara = p;\n  } else if (t.indexOf(\"Randomize\") === 0 && t.indexOf(\"blocks\") !== -1) {\n    blocksPara = p;\n  }\n}\n\n// 1) \"Add custom field size\" -> prepend a \"Done: \" run.\nif (addCustomFieldPara) {\n  const startRange = addCustomFieldPara.getRange(Word.RangeLocation.start);\n  startRange.insertOoxml(doneRunOoxml(\"Done: \"), Word.InsertLocation.before);\n}\n\n// 2) \"Randomize \u2013 blocks / num of blocks\" -> prepend a \"Done: \" run\n//    (existing runs / proofErr spell-check markers stay untouched).\nif (blocksPara) {\n  const startRange = blocksPara.getRange(Word.RangeLocation.start);\n  startRange.insertOoxml(doneRunOoxml(\"Done: \"), Word.InsertLocation.before);\n}\n\nawait context.sync();\n\n// 3) \"Randomize \u2013 zombie and human\" -> Word stamps the cursor's last\n//    edit location with the implicit \"_GoBack\" bookmark when the\n//    document is saved after an edit.\nif (zombieHumanPara) {\n  const bmRange = zombieHumanPara.getRange(Word.RangeLocation.start);\n  bmRange.insertBookmark(\"_GoBack\");\n}\n\nawait context.sync();\n", "ps1": "# Mark two maintenance-request items as done, and add \"Done: \" prefixes:\n#   - \"Add custom field size\"               -> \"Done: Add custom field size\"\n#   - \"Randomize \u2013 blocks / num of blocks\"  -> \"Done: Randomize \u2013 blocks / num of blocks\"\n# (matches commit message: \"Made change to randomize both the human and zombie\")\n\n$d = $word.ActiveDocument\n\nfunction Find-Paragraph($doc, [string]$exactText) {\n    foreach ($p in $doc.Paragraphs) {\n        $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($t -eq $exactText) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# 1) \"Add custom field size\" -> prepend a \"Done: \" run (kept as its OWN run,\n#    not merged into the existing \"Add custom field size\" run). Replacing\n#    the whole paragraph range (including its end-of-paragraph mark) with\n#    OOXML that contains two separate <w:r> elements keeps them distinct.\n$pAdd = Find-Paragraph $d \"Add custom field size\"\nif ($pAdd -ne $null) {\n    $rng = $pAdd.Range\n    $ooxml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">Done: </w:t></w:r><w:r><w:t>Add custom field size</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $rng.InsertXML($ooxml)\n}\n\n# 2) \"Randomize \u2013 blocks / num of blocks\" -> prepend a \"Done: \" run, keeping\n#    the rest of the paragraph (the \"num\" spell-check-flagged run and its\n#    surrounding proofErr markers) exactly as it was.\n$pBlocks = Find-Paragraph $d \"Randomize \u2013 blocks / num of blocks\"\nif ($pBlocks -ne $null) {\n    $rng = $pBlocks.Range\n    $ooxml = '<?xml version=\"1.0\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t xml:space=\"preserve\">Done: </w:t></w:r><w:r><w:t xml:space=\"preserve\">Randomize \u2013 blocks / </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>num</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> of blocks</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n    $rng.InsertXML($ooxml)\n}\n\n# 3) \"Randomize \u2013 zombie and human\" -> Word stamps the cursor's last edit\n#    location with the implicit \"_GoBack\" bookmark when the document is\n#    saved after an edit.\n$pZombie = Find-Paragraph $d \"Randomize \u2013 zombie and human\"\nif ($pZombie -ne $null) {\n    $bmRng = $pZombie.Range\n    $bmRng.Collapse(1)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRng)\n}\n"}
